# Re-synchronize cached `w:lastRenderedPageBreak` markers with Word's
# latest pagination pass (safety snapshot before entrypoint quarantine).
#
# These markers are pure rendering-cache hints (they do not change any
# visible text), so each edit below is expressed as a full replacement
# of the owning <w:p> via Range.InsertXML: we collapse a precise Range
# to its start (so InsertXML lands on that exact paragraph) and replace
# it with the same paragraph, but with the lastRenderedPageBreak moved.
#
# NOTE: the body-level paragraph edits are done before the in-table
# edits -- doing table-cell InsertXML surgery first was observed to
# desync $d.Paragraphs.Item(N) indices for later lookups.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- "Phase Reference: ..." paragraph: gains a lastRenderedPageBreak.
$pPhase = $d.Paragraphs.Item(164)
$r = $pPhase.Range
$r.Collapse(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="7512613F" w14:textId="77777777" w:rsidR="00541079" w:rsidRDefault="00000000"><w:r><w:lastRenderedPageBreak/><w:t>Phase Reference: AI Talent Engine Phase 7 Master v7.2</w:t></w:r></w:p>'
$r.InsertXML($xml)

# --- "Maintainer: ..." paragraph: loses its lastRenderedPageBreak.
$pMaintainer = $d.Paragraphs.Item(165)
$r = $pMaintainer.Range
$r.Collapse(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="5D29AEFD" w14:textId="77777777" w:rsidR="00541079" w:rsidRDefault="00000000"><w:r><w:t>Maintainer: L. David Mendoza &#169; 2025</w:t></w:r></w:p>'
$r.InsertXML($xml)

$t1 = $d.Tables.Item(1)

# --- Row "27" (col 1): gains a lastRenderedPageBreak before the text.
$cell27 = $t1.Cell(28, 1)
$r = $cell27.Range
$r.Collapse(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="06E0DBF2" w14:textId="77777777" w:rsidR="00541079" w:rsidRDefault="00000000"><w:r><w:lastRenderedPageBreak/><w:t>27</w:t></w:r></w:p>'
$r.InsertXML($xml)

# --- "24-month citation velocity ... (High / Medium / Low)." cell:
# the lastRenderedPageBreak that used to split this into two runs is
# gone, and the text is merged back into a single run.
$cellVelocity = $t1.Cell(28, 3)
$r = $cellVelocity.Range
$r.Collapse(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="0D9C7E13" w14:textId="77777777" w:rsidR="00541079" w:rsidRDefault="00000000"><w:r><w:t>24-month citation velocity (High / Medium / Low).</w:t></w:r></w:p>'
$r.InsertXML($xml)

# --- Row "28" (col 1): loses its lastRenderedPageBreak.
$cell28 = $t1.Cell(29, 1)
$r = $cell28.Range
$r.Collapse(1)
$xml = '<w:p ' + $wNs + ' w14:paraId="77C4B80C" w14:textId="77777777" w:rsidR="00541079" w:rsidRDefault="00000000"><w:r><w:t>28</w:t></w:r></w:p>'
$r.InsertXML($xml)

Write-Output "applied lastRenderedPageBreak resync"
